# Apply the requested edits to the "Frameworks and Libraries" / "Development
# Tools" bullet items: collapse runs that were split around spell-checked
# words (TensorFlow/PyTorch, Matplotlib/Seaborn/Plotly, Jupyter) back into a
# single run each, and drop the "or Flask" alternative from the Backend line.

$d = $word.ActiveDocument

# 1) Deep Learning: "TensorFlow or " + "PyTorch" (spell-checked) + " for model
#    development and training." -> single run, no proofErr wrapper.
$d.Content.Find.Execute(
    "TensorFlow or PyTorch for model development and training.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "TensorFlow or PyTorch for model development and training.", 2)

# 2) Backend: remove the "or Flask" alternative.
$d.Content.Find.Execute(
    ": Django or Flask for building RESTful APIs.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Django for building RESTful APIs.", 2)

# 3) Data Visualization: "Matplotlib, Seaborn, or " + "Plotly" (spell-checked)
#    + " for presenting model results and metrics graphically." -> single run.
$d.Content.Find.Execute(
    "Matplotlib, Seaborn, or Plotly for presenting model results and metrics graphically.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Matplotlib, Seaborn, or Plotly for presenting model results and metrics graphically.", 2)

# 4) Integrated Development Environment (IDE): ": VS Code, PyCharm, or " +
#    "Jupyter" (spell-checked) + " Notebooks for coding and testing." ->
#    single run.
$d.Content.Find.Execute(
    ": VS Code, PyCharm, or Jupyter Notebooks for coding and testing.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": VS Code, PyCharm, or Jupyter Notebooks for coding and testing.", 2)
